$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
if (-not $ws) { $ws = $wb.ActiveSheet }

$row = 44

# Force column A to be treated as plain text so the date-like string
# "10/15/2025" is stored verbatim instead of being auto-converted into a
# date serial number, matching the existing rows above it. Reset the
# style back to "Normal" afterwards so no stray number-format style gets
# attached to the new cell (keeping it identical in shape to the other
# data rows, which carry no explicit style index).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/15/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1768812452900025
$ws.Cells.Item($row, 3).Value = 0.8231187547099975
